$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "82.0-100.0"

# Row 3
$ws.Range("C3").Value = "27.0-34.0"

# Row 4
$ws.Range("C4").Value = "316-354"

# Row 5
$ws.Range("A5").Value = "RBC分布宽度"
$ws.Range("C5").Value = "37.0-50.0"

# Row 6
$ws.Range("A6").Value = "RBC分布宽度"
$ws.Range("C6").Value = "11.6-14.8"

# Row 7
$ws.Range("A7").Value = "血小板计数"
$ws.Range("C7").Value = "125-350"

# Row 8
$ws.Range("A8").Value = "血小板分布宽度"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "11.2"
$ws.Range("C8").Value = "9.00-17.00"

# Row 9
$ws.Range("C9").Value = "9-13"
$ws.Range("D9").Value = "fL"

# Row 10
$ws.Range("C10").Value = "0.17-0.35"

# Row 11
$ws.Range("A11").Value = "大血小板比值"
$ws.Range("C11").Value = "15-45"
